$d = $word.ActiveDocument

function Replace-InRange($range, $oldText, $newText) {
    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# ----------------------------------------------------------------------
# Title
# ----------------------------------------------------------------------
Replace-InRange $d.Paragraphs(1).Range `
    "Sacred Texts in Ancient Rome: Facilitators of Cohesion and Disintegration" `
    "Democracy: Ensuring Equal Voices and Shared Power"

# ----------------------------------------------------------------------
# Author name: "Dr" + "." + " John Stevens" (3 runs) -> "Claire Armstrong" (1 run)
# ----------------------------------------------------------------------
$r = $d.Paragraphs(2).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Text = "Claire Armstrong"

# ----------------------------------------------------------------------
# Email: "jstevens@classicstudies" + "." + "edu" -> "claire" + "." + "armstrong@edumail" + "." + "org"
# ----------------------------------------------------------------------
Replace-InRange $d.Paragraphs(3).Range "jstevens@classicstudies" "claire"
Replace-InRange $d.Paragraphs(3).Range "edu" "armstrong@edumail"

$r = $d.Paragraphs(3).Range
$r.MoveEnd(1, -1) | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter(".") | Out-Null
$r.Collapse(0) | Out-Null
$r.InsertAfter("org") | Out-Null

# ----------------------------------------------------------------------
# Body paragraph - first 4 "sentences" before the first manual line-break pair
# ----------------------------------------------------------------------
$body = $d.Paragraphs(5).Range

Replace-InRange $body `
    "In the intricate tapestry of ancient Roman society, sacred texts served as both the threads that wove the fabric of unity and the scissors that severed it" `
    "In the realm of governance, where authority and decision-making are exercised, the concept of democracy takes center stage"

Replace-InRange $body `
    " They were the conduits of tradition, morality, and religious fervor, weaving a collective identity among the diverse populace" `
    " As a form of government, democracy empowers citizens with the right to participate in the selection of their leaders and policies that govern their lives"

Replace-InRange $body `
    " Through their revered scriptures, Romans found solace, guidance, and a sense of belonging in a world perpetually teetering on the precipice of chaos" `
    " Whether ancient Athenian assemblies or modern-day representative democracies, this system has been consistently advocated for its commitment to promoting equality, liberty, and shared power"

Replace-InRange $body `
    " Yet, these same texts could foment discord and disintegration, wielding the power to ignite sectarian strife and political machinations" `
    " In essence, democracy upholds the belief that the most effective way to govern a society is through the active participation and consent of its citizens"

Replace-InRange $body `
    "In the annals of Roman history, the venerable Sibylline Books stand as illuminating testaments to the dual capacity of sacred texts to unite and divide" `
    "Democracy, in its true essence, is characterized by the fundamental principles of equality, inclusivity, and shared power"

Replace-InRange $body `
    " These enigmatic tomes, imbued with oracular pronouncements of impending doom and divine favor, were consulted by the Senate in times of crisis, their cryptic verses interpreted to discern the will of the gods" `
    " It is a system that recognizes the inherent worth and dignity of every individual, regardless of their race, gender, ethnicity, or creed"

Replace-InRange $body `
    " The curators of these sacred texts, the quindecimviri sacris faciundis, were not merely custodians of ancient lore but custodians of the city's destiny" `
    " Through democratic processes, citizens are empowered to make collective decisions, holding elected officials accountable for their actions"

Replace-InRange $body `
    " Their pronouncements had the power to rally the populace, galvanize legions, and steer the course of history" `
    " By embracing diversity and fostering inclusivity, democracy cultivates a sense of unity and belonging among citizens, allowing them to work together towards the common good"

Replace-InRange $body `
    "Furthermore, the advent of Christianity in the Roman Empire introduced a new dimension to the complex relationship between sacred texts and societal cohesion" `
    "Moreover, the strength of a democracy lies in the active participation of its citizens"

Replace-InRange $body `
    " The rise of the Christian faith, with its revolutionary teachings of pacifism and universal love, sparked a profound clash with the prevailing Roman ethos of martial prowess and civic duty" `
    " It is not merely a spectator sport where citizens passively observe the actions of their leaders"

Replace-InRange $body `
    " The incompatibility of these belief systems ignited fierce persecutions, tearing at the seams of Roman society and culminating in the infamous Great Persecution under Emperor Diocletian, where sacred texts became both weapons of oppression and instruments of unwavering faith" `
    " Instead, it entails active engagement, dialogue, and deliberation"

# ----------------------------------------------------------------------
# Insert the large new block of sentences right before the paragraph's
# trailing "." run (which stays as-is).
# ----------------------------------------------------------------------
$body = $d.Paragraphs(5).Range
$insertPoint = $body.Duplicate
$insertPoint.MoveEnd(1, -1) | Out-Null   # drop the paragraph mark
$insertPoint.Collapse(0) | Out-Null      # collapse to just before the very end
$insertPoint.MoveStart(1, -1) | Out-Null # move start back one char so this range covers the trailing "."
$insertPoint.Collapse(1) | Out-Null      # collapse to just before that trailing "."

$insertPoint.InsertAfter(" Democracy provides a platform for citizens to voice their opinions, hold their representatives accountable, and shape the direction of their society")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" By encouraging civic engagement and participation, democracy empowers individuals to influence decisions that directly impact their lives, fostering a sense of ownership and responsibility")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter([char]11)
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter([char]11)
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter("Body:")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter([char]11)
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter("Exploring the historical evolution of democracy reveals a fascinating narrative of struggle, innovation, and adaptation")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" From the ancient Greek city-states, where direct democracy allowed citizens to participate directly in decision-making, to the representative democracies that emerged in response to growing populations and geographical challenges, democracy has undergone significant transformations")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" The American Revolution, the French Revolution, and the countless movements for independence and self-governance around the world stand as testament to the enduring power of the democratic ideal")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter([char]11)
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter([char]11)
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter("The functioning of a democracy is intricate and multifaceted")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" At the core of many democracies lies the concept of separation of powers, ensuring that no single branch of government holds absolute authority")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" Executive, legislative, and judicial branches, each with distinct roles and responsibilities, work in concert to maintain checks and balances")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" Additionally, the presence of independent institutions such as a free press and a robust civil society plays a crucial role in holding governments accountable, safeguarding individual rights, and ensuring transparency")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter([char]11)
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter([char]11)
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter("While democracy offers a framework for just and equitable governance, it is not without its challenges")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" Ensuring fair representation, addressing the influence of money in politics, promoting civic engagement, and combating disinformation are just a few of the hurdles that democracies face")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(".")
$insertPoint.Collapse(0) | Out-Null
$insertPoint.InsertAfter(" Additionally, the rise of populism and authoritarianism in recent years has posed a significant threat to the principles of democratic governance, emphasizing the need for vigilance and unwavering commitment to democratic values")

# ----------------------------------------------------------------------
# Summary heading paragraph stays the same ("Summary") - no change.
# Summary body paragraph.
# ----------------------------------------------------------------------
$summary = $d.Paragraphs(7).Range

Replace-InRange $summary `
    "Sacred texts played a paradoxical role in ancient Rome, serving as both architects of unity and agents of discord" `
    "Democracy, as a form of government, stands as a testament to the power of collective decision-making and the inherent worth of each individual"

Replace-InRange $summary `
    " The Sibylline Books exemplified this duality, embodying the power to galvanize the populace and shape the city's destiny while simultaneously harboring the potential for political manipulation" `
    " Through the principles of equality, inclusivity, and shared power, democracy empowers citizens to shape their societies and hold their leaders accountable"

Replace-InRange $summary `
    " The introduction of Christianity further exacerbated this dichotomy, instigating a clash between entrenched traditions and revolutionary ideals" `
    " Throughout history, democracy has undergone profound transformations, adapting to changing circumstances while remaining steadfast in its commitment to liberty, justice, and the pursuit of the common good"

Replace-InRange $summary `
    " This profound clash manifested in vehement persecutions, highlighting the explosive force of sacred texts when wielded as instruments of power and belief" `
    " While challenges remain, the enduring strength of democracy lies in its ability to embrace diversity, foster civic engagement, and inspire citizens to work together for a better future"

# Remove the now-obsolete trailing sentence + its period run.
Replace-InRange $summary " Thus, in ancient Rome, sacred texts were both the threads that bound and the wedges that sundered" ""

$summaryParas = $d.Paragraphs
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastText = $lastPara.Range.Text
if ($lastText.Trim().Length -eq 0) {
    # already has a trailing empty paragraph (shouldn't normally happen)
} else {
    $lastPara.Range.InsertParagraphAfter() | Out-Null
}

Write-Output "done"
